# Generate Report for Handback
# Adds a new handed-back file (5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md) as
# row 3 to the Overview / zh-cn / de-de sheets, extending their tables and
# hyperlinks to match.

$wb = $excel.ActiveWorkbook

$newFile      = "5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md"
$newPath      = "e2e\5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md"
$newStatus    = "Handed back: in sync with en-US"
$dateFormat   = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Range("A3").Value = $newFile
$ws1.Range("B3").Value = $newPath
$ws1.Range("C3").Value = ".md"
$ws1.Range("E3").Value = $newStatus
$ws1.Range("F3").Value = $newStatus
$ws1.Range("G3").Value = "2016-10-26 07:21:28"
$ws1.Range("G3").NumberFormat = $dateFormat

$ws1.Hyperlinks.Add($ws1.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f376d12adacbe7b0cf85bd3de2244ff0860a9a53/e2e/5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md", "", "", $newPath)

$ws1.ListObjects.Item(1).Resize($ws1.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Range("A3").Value = $newFile
$ws2.Range("B3").Value = ".md"
$ws2.Range("C3").Value = $newStatus
$ws2.Range("D3").Value = "e2e"
$ws2.Range("E3").Value = "ht"
$ws2.Range("F3").Value = "'True"
$ws2.Range("G3").Value = "5e7c082f-c6cf-4944-831d-a90fb88d9e3c.1b8e8a361b8180533ba7730e2eaefae6cfd9d0c9.zh-cn.xlf"
$ws2.Range("H3").Value = "2016-10-26 07:21:16"
$ws2.Range("H3").NumberFormat = $dateFormat
$ws2.Range("I3").Value = $newFile
$ws2.Range("J3").Value = "5e7c082f-c6cf-4944-831d-a90fb88d9e3c.1b8e8a361b8180533ba7730e2eaefae6cfd9d0c9.zh-cn.xlf"
$ws2.Range("K3").Value = "2016-10-26 07:21:57"
$ws2.Range("K3").NumberFormat = $dateFormat
$ws2.Range("L3").Value = "'"
$ws2.Range("M3").Value = "'True"
$ws2.Range("N3").Value = "'"
$ws2.Range("O3").Value = "'False"
$ws2.Range("P3").Value = "'"

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f376d12adacbe7b0cf85bd3de2244ff0860a9a53/e2e/5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md", "", "", $newFile)
$ws2.Hyperlinks.Add($ws2.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/8dd634329835c70953255d53cada1a846189f887/e2e/5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md", "", "", $newFile)

$ws2.ListObjects.Item(1).Resize($ws2.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Range("A3").Value = $newFile
$ws3.Range("B3").Value = ".md"
$ws3.Range("C3").Value = $newStatus
$ws3.Range("D3").Value = "e2e"
$ws3.Range("E3").Value = "ht"
$ws3.Range("F3").Value = "'True"
$ws3.Range("G3").Value = "5e7c082f-c6cf-4944-831d-a90fb88d9e3c.1b8e8a361b8180533ba7730e2eaefae6cfd9d0c9.de-de.xlf"
$ws3.Range("H3").Value = "2016-10-26 07:21:28"
$ws3.Range("H3").NumberFormat = $dateFormat
$ws3.Range("I3").Value = $newFile
$ws3.Range("J3").Value = "5e7c082f-c6cf-4944-831d-a90fb88d9e3c.1b8e8a361b8180533ba7730e2eaefae6cfd9d0c9.de-de.xlf"
$ws3.Range("K3").Value = "2016-10-26 07:22:15"
$ws3.Range("K3").NumberFormat = $dateFormat
$ws3.Range("L3").Value = "'"
$ws3.Range("M3").Value = "'True"
$ws3.Range("N3").Value = "'"
$ws3.Range("O3").Value = "'False"
$ws3.Range("P3").Value = "'"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f376d12adacbe7b0cf85bd3de2244ff0860a9a53/e2e/5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md", "", "", $newFile)
$ws3.Hyperlinks.Add($ws3.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/55bd1f2994a94b3ed53739ec7f331c77749a9a55/e2e/5e7c082f-c6cf-4944-831d-a90fb88d9e3c.md", "", "", $newFile)

$ws3.ListObjects.Item(1).Resize($ws3.Range("A1:P3"))

Write-Host "Handback report row added to Overview, zh-cn and de-de sheets."
